# Updated symbol list on Tue Feb 14 18:34:52 UTC 2023 with GitHub Actions
# Applies updated Price (D) and Volume(1h) (E) values for the crypto table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'296.49"
$ws.Range("E2").Value = "'2.80%"
$ws.Range("D3").Value = "'41.48"
$ws.Range("E3").Value = "'3.00%"
$ws.Range("D4").Value = "'5.039"
$ws.Range("E4").Value = "'0.04%"
$ws.Range("D5").Value = "'0.07496"
$ws.Range("E5").Value = "'2.83%"
$ws.Range("E6").Value = "'4.25%"
$ws.Range("D7").Value = "'0.9262"
$ws.Range("E7").Value = "'1.16%"
$ws.Range("D9").Value = "'0.1219"
$ws.Range("E9").Value = "'2.41%"
$ws.Range("D10").Value = "'0.1843"
$ws.Range("E10").Value = "'7.80%"
$ws.Range("D11").Value = "'0.08953"
$ws.Range("E11").Value = "'3.86%"
$ws.Range("D12").Value = "'0.04144"
$ws.Range("E12").Value = "'-0.69%"
$ws.Range("D13").Value = "'0.1053"
$ws.Range("E13").Value = "'-0.10%"
$ws.Range("D14").Value = "'0.001287"
$ws.Range("E14").Value = "'0.41%"
$ws.Range("D15").Value = "'0.005859"
$ws.Range("E15").Value = "'-0.47%"
$ws.Range("E16").Value = "'-1.68%"
$ws.Range("D17").Value = "'4.362"
$ws.Range("E17").Value = "'1.91%"
$ws.Range("D18").Value = "'0.3317"
$ws.Range("E18").Value = "'1.70%"
$ws.Range("D19").Value = "'7.945"
$ws.Range("E19").Value = "'1.99%"
$ws.Range("E20").Value = "'4.19%"
$ws.Range("D21").Value = "'0.2964"
$ws.Range("E21").Value = "'2.75%"
$ws.Range("D22").Value = "'0.04041"
$ws.Range("E22").Value = "'4.85%"
$ws.Range("D23").Value = "'0.001267"
$ws.Range("E23").Value = "'-0.04%"
$ws.Range("D24").Value = "'0.003879"
$ws.Range("E24").Value = "'2.13%"
$ws.Range("D25").Value = "'0.0001229"
$ws.Range("E25").Value = "'-4.01%"
$ws.Range("E26").Value = "'0.00%"
$ws.Range("D38").Value = "'0.02415"
$ws.Range("E38").Value = "'4.96%"
$ws.Range("D39").Value = "'0.05208"
$ws.Range("E39").Value = "'5.26%"
$ws.Range("D40").Value = "'0.006407"
$ws.Range("E40").Value = "'-9.78%"
$ws.Range("D41").Value = "'0.007793"
$ws.Range("E41").Value = "'1.24%"
$ws.Range("D42").Value = "'0.1325"
$ws.Range("E42").Value = "'4.40%"
$ws.Range("D43").Value = "'0.007369"
$ws.Range("E43").Value = "'0.05%"
$ws.Range("D44").Value = "'0.008107"
$ws.Range("E44").Value = "'6.62%"
$ws.Range("D45").Value = "'0.2973"
$ws.Range("E45").Value = "'-4.67%"
$ws.Range("D46").Value = "'0.00006254"
$ws.Range("E46").Value = "'-1.76%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.11%"
$ws.Range("D48").Value = "'0.04556"
$ws.Range("E48").Value = "'-81.13%"
$ws.Range("D49").Value = "'0.004198"
$ws.Range("E49").Value = "'-0.09%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.11%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.11%"

Write-Host "Applied cryptos price/volume updates."
